$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers corresponding to B,C,D,E,F,I,J,K,L,M,N (skipping G and H which are unchanged/empty)
$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

$data = @(
    @(1.02, 1.030467032709939, 1.040428727334756, 1.030144171202675, 1.041724324102612, 1.039609291165755, 1.035608117816501, 1.043211102458997, 1.03295603444136, 1.044503029288581, 1.01577726647918),
    @(1.02, 1.031626290725917, 1.041115305036162, 1.031135725631961, 1.043086943986686, 1.03994444868674, 1.036407771802833, 1.043708427514884, 1.033755389573087, 1.045674887641358, 1.016049898131116),
    @(1.02, 1.032375996687559, 1.04155910439444, 1.031777287015109, 1.043968488849291, 1.040159710106165, 1.036924282966968, 1.044029091502155, 1.034271977635892, 1.046432467684437, 1.016225808845747),
    @(1.02, 1.032691076535835, 1.041745566636252, 1.032046990657591, 1.044339054587261, 1.040249821401146, 1.037141205511406, 1.044163626073133, 1.034488997057129, 1.046750791435959, 1.016299642241884),
    @(1.02, 1.032743974225064, 1.041776867947872, 1.032092274600031, 1.044401272194598, 1.040264928940159, 1.037177614966243, 1.044186199036126, 1.034525426558629, 1.046804229913184, 1.016312032196393),
    @(1.02, 1.032380207176777, 1.041561596348659, 1.031780890841614, 1.043973440503571, 1.040160915688533, 1.036927182354761, 1.044030890230693, 1.03427487806472, 1.046436721779017, 1.016226795880282),
    @(1.02, 1.030858895775713, 1.040660855017296, 1.030479280661227, 1.042184863019906, 1.039722892347565, 1.035878555473159, 1.043379411548084, 1.033226315041709, 1.044899208789562, 1.015869507359673),
    @(1.02, 1.02817491038172, 1.039070106449038, 1.028185322168162, 1.039031795394418, 1.03893870933315, 1.034023655852822, 1.042222701563088, 1.031373604954911, 1.042184514681324, 1.01523607572134),
    @(1.02, 1.026383281595701, 1.038007251451614, 1.026655713775298, 1.036928648047903, 1.03840760761105, 1.032782218959062, 1.041445693649502, 1.030135034001505, 1.040370931950171, 1.014811184983073),
    @(1.02, 1.025606910867721, 1.03754646710621, 1.025993289557331, 1.036017658911312, 1.038175657257951, 1.032243498683045, 1.041107847102456, 1.029597890432748, 1.039584697003364, 1.014626580733951),
    @(1.02, 1.025318441935231, 1.037375227074988, 1.025747220234151, 1.03567922635686, 1.038089202556266, 1.03204321665387, 1.040982145721429, 1.029398244694715, 1.039292509711537, 1.014557916358295),
    @(1.02, 1.025380323627373, 1.03741196245223, 1.025800003662062, 1.035751823638921, 1.038107760878175, 1.032086185878895, 1.041009118607299, 1.029441075138855, 1.039355191472838, 1.01457264936677),
    @(1.02, 1.025583067792397, 1.037532314069148, 1.025972949714102, 1.035989685024924, 1.038168516968028, 1.032226946934645, 1.041097460875157, 1.02958139024639, 1.039560547679262, 1.014620906838098),
    @(1.02, 1.025707973123744, 1.03760645551653, 1.026079505392251, 1.036136232456213, 1.038205911293449, 1.03231365092157, 1.041151863610108, 1.029667826184846, 1.039687055156878, 1.014650627371899),
    @(1.02, 1.026434794246434, 1.03803782037711, 1.026699674595903, 1.036989100499354, 1.03842295963137, 1.03281794726738, 1.041468085958096, 1.030170664758847, 1.040423091665779, 1.014823423381008),
    @(1.02, 1.02689055185184, 1.038308253999517, 1.027088664285311, 1.037523995841074, 1.038558577869463, 1.033133964721094, 1.041666069695175, 1.030485857811756, 1.040884533593134, 1.014931646455626),
    @(1.02, 1.027156331413722, 1.038465939210763, 1.027315546444535, 1.037835961344495, 1.038637490709524, 1.033318179581554, 1.04178141548694, 1.030669624193988, 1.041153594218291, 1.014994710989406),
    @(1.02, 1.027246945932555, 1.038519696605874, 1.027392905922767, 1.037942328559609, 1.038664365571802, 1.033380973007587, 1.041820722568819, 1.03073227018707, 1.041245321659218, 1.015016204175328),
    @(1.02, 1.026841659199078, 1.038279244630891, 1.027046930308128, 1.037466609751942, 1.03854404707433, 1.033100070714859, 1.041644841855695, 1.030452048905263, 1.040835034639491, 1.01492004136901),
    @(1.02, 1.025523367201982, 1.037496875841539, 1.025922021850781, 1.035919642202503, 1.038150634045091, 1.032185501219961, 1.041071452081821, 1.029540074453114, 1.039500079412586, 1.014606698817681),
    @(1.02, 1.024693981100814, 1.037004482362875, 1.025214657028684, 1.03494670755728, 1.037901555162246, 1.031609448308867, 1.040709723292267, 1.028965946183836, 1.038659901749458, 1.014409143052806),
    @(1.02, 1.025133704864581, 1.037265555612202, 1.025589653526571, 1.035462507794081, 1.038033760244877, 1.031914922707706, 1.04090159790433, 1.029270372453221, 1.039105376345238, 1.014513922891613),
    @(1.02, 1.026863751857093, 1.03829235288306, 1.027065788123223, 1.037492540146941, 1.038550613505083, 1.033115386305334, 1.041654434218874, 1.030467325942819, 1.040857401365261, 1.014925285393851),
    @(1.02, 1.028869181941209, 1.039481769436754, 1.028778414061988, 1.03984711970221, 1.039142902958222, 1.034504039029463, 1.042522772743085, 1.031853174942488, 1.042886984123034, 1.015400290392737)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $cols[$c]).Value = $rowVals[$c]
    }
}

Write-Host "Updated vm_pu values for 380 kV case"
